$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-record entries (rows 122-123): "Adding new puzzle" / "Fixing buggs of new puzzle" ---
$ws.Range("A122").Value = "2022-02-25"
$ws.Range("B122").Value = "Adding new puzzle"
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = "Sarvan Amel"

$ws.Range("A123").Value = "2022-02-26"
$ws.Range("B123").Value = "Fixing buggs of new puzzle"
$ws.Range("C123").Value = 2
$ws.Range("D123").Value = "Sarvan Amel"

# --- Update workhours total for Sarvan Amel (row 129) and add a C129 number cell ---
$ws.Range("B129").Value = 13.58
$ws.Range("C129").NumberFormat = "0.00"

# --- Row 130 (grand total row): thin out the bottom border on C130 and drop the ---
# --- explicit "thick bottom" row formatting so it matches the rest of the table ---
$ws.Range("C130").Borders.Item(9).LineStyle = 1
$ws.Range("C130").Borders.Item(9).Weight = 2
$ws.Rows.Item(130).AutoFit()

# --- Move the active view/selection further down the sheet ---
$null = $ws.Range("G117").Select()
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 1
